$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date placeholder text (Master +
#    every Custom Layout) from 30/12/2017 -> 14/4/2018, same as what real
#    PowerPoint does to the cached field text when the deck is re-saved on a
#    later date.
# ---------------------------------------------------------------------------
$newDate = "14/4/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -ne $newDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DatePlaceholder $layouts.Item($l).Shapes
}

# ---------------------------------------------------------------------------
# 2) Rename the "prevAddressBook" variable shown in the three command-object
#    tables on slide 1 to "prevImdb" (keeping the trailing " = sN" text).
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellTr = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
                for ($pi = 1; $pi -le $cellTr.Paragraphs().Count; $pi++) {
                    $para = $cellTr.Paragraphs($pi)
                    if ($para.Text -match "^prevAddressBook( = s\d+)$") {
                        $para.Text = "prevImdb" + $matches[1]
                    }
                }
            }
        }
    }
}
